$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New observation row for 2023-05-01 (Boucherville, BUAM frog chorus).
# Row 36 already carries the exact style/formatting the new row needs,
# so copy it (values + formats) into row 37, then overwrite the values.
$ws.Range("A36:I36").Copy($ws.Range("A37:I37"))

$ws.Cells.Item(37, 1).Value = 45047
$ws.Cells.Item(37, 2).Value = "BUAM"
$ws.Cells.Item(37, 3).Value = "N/A"
$ws.Cells.Item(37, 4).Value = "Boucherville"
$ws.Cells.Item(37, 5).Value = "Montérégie"
$ws.Cells.Item(37, 6).Value = "A"
$ws.Cells.Item(37, 7).Value = "Cote 3"
$ws.Cells.Item(37, 9).Value = "Ornitholarocque"
$ws.Cells.Item(37, 8).Value = "iNaturalist (https://www.inaturalist.org/observations/159092353)"

$ws.Range("E43").Select()
